$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 118, shifting the existing rows 118-129
# down to 120-131 (and carrying their formatting, e.g. the date style on
# column D, down with them).
$ws.Rows("118:119").Insert()

# Row 118: new "Camote" record dated 2021-09-10 (serial 44449)
$ws.Range("A118").Value = 5
$ws.Range("B118").Value = "Macroferia Regional de Talca"
$ws.Range("C118").Value = "Maule"
$ws.Range("D118").Value = 44449
$ws.Range("E118").Value = 7
$ws.Range("F118").Value = 100112045
$ws.Range("G118").Value = "Zapallo"
$ws.Range("H118").Value = "Camote"
$ws.Range("I118").Value = "1a (guarda)"
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 600
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 600
$ws.Range("N118").Value = "$/kilo (volumen en unidades)"
$ws.Range("O118").Value = "Región del Maule"
$ws.Range("P118").Value = 600
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = "Hortaliza"

# Row 119: new "Paine" record dated 2021-09-10 (serial 44449)
$ws.Range("A119").Value = 5
$ws.Range("B119").Value = "Macroferia Regional de Talca"
$ws.Range("C119").Value = "Maule"
$ws.Range("D119").Value = 44449
$ws.Range("E119").Value = 7
$ws.Range("F119").Value = 100112045
$ws.Range("G119").Value = "Zapallo"
$ws.Range("H119").Value = "Paine"
$ws.Range("I119").Value = "1a (guarda)"
$ws.Range("J119").Value = 1200
$ws.Range("K119").Value = 140
$ws.Range("L119").Value = 140
$ws.Range("M119").Value = 140
$ws.Range("N119").Value = "$/kilo (volumen en unidades)"
$ws.Range("O119").Value = "Región del Maule"
$ws.Range("P119").Value = 140
$ws.Range("Q119").Value = 1
$ws.Range("R119").Value = "Hortaliza"
